$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 25: date value, reuse the date number-format style already used by A1
$ws.Range("A1").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A25").Value = 44930

# Row 26: time value (reuse time style from A3) + description
$ws.Range("A3").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("A26").Value = 0.5625
$ws.Range("B26").Value = "SqlDataReader 原理 只傳Data索引值"

# Row 27: time value + description
$ws.Range("A3").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A27").Value = 0.61805555555555558
$ws.Range("B27").Value = "SQL語法 切記不要組合字串"

# Row 28: time value + description
$ws.Range("A3").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A28").Value = 0.625
$ws.Range("B28").Value = "登入程式 寫法"

$excel.CutCopyMode = 0

# Update the view/selection to match the author's final state
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("I26").Select()
